$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Cells.Item(51, 8).Value = 4393.75
$ws.Cells.Item(51, 9).Value = 3040
$ws.Cells.Item(51, 10).Value = 5747.5
$ws.Cells.Item(51, 11).Value = 3040
$ws.Cells.Item(51, 12).Value = 5747.5
$ws.Cells.Item(51, 13).Value = -2556
$ws.Cells.Item(51, 14).Value = -6715.5
# Row 52
$ws.Cells.Item(52, 8).Value = 249.08824
$ws.Cells.Item(52, 9).Value = 59.333332
$ws.Cells.Item(52, 10).Value = 289.75
$ws.Cells.Item(52, 11).Value = 177.999996
$ws.Cells.Item(52, 12).Value = 869.25
$ws.Cells.Item(52, 13).Value = -17.99999600000001
$ws.Cells.Item(52, 14).Value = -1189.25
# Row 98
$ws.Cells.Item(98, 8).Value = 2844.4348
$ws.Cells.Item(98, 9).Value = 3087
$ws.Cells.Item(98, 11).Value = 3087
$ws.Cells.Item(98, 13).Value = -1589
# Row 100
$ws.Cells.Item(100, 8).Value = 2142.923
$ws.Cells.Item(100, 9).Value = 2293.75
$ws.Cells.Item(100, 11).Value = 2293.75
$ws.Cells.Item(100, 13).Value = -1752.75
# Row 122
$ws.Cells.Item(122, 8).Value = 2844.4348
$ws.Cells.Item(122, 9).Value = 3087
$ws.Cells.Item(122, 11).Value = 9261
$ws.Cells.Item(122, 13).Value = -6811
# Row 132
$ws.Cells.Item(132, 8).Value = 1698.8334
$ws.Cells.Item(132, 9).Value = 1254.4
$ws.Cells.Item(132, 11).Value = 3763.2
$ws.Cells.Item(132, 13).Value = -1233.2
# Row 138
$ws.Cells.Item(138, 8).Value = 1728527
$ws.Cells.Item(138, 9).Value = 2654.6
$ws.Cells.Item(138, 10).Value = 2636881
$ws.Cells.Item(138, 11).Value = 7963.799999999999
$ws.Cells.Item(138, 12).Value = 7910643
$ws.Cells.Item(138, 13).Value = -2823.799999999999
$ws.Cells.Item(138, 14).Value = -7920923

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 5104.857
$ws.Cells.Item(2, 9).Value = 7003
$ws.Cells.Item(2, 10).Value = 3681.25
$ws.Cells.Item(2, 11).Value = 7003
$ws.Cells.Item(2, 12).Value = 3681.25
$ws.Cells.Item(2, 13).Value = -6890
$ws.Cells.Item(2, 14).Value = -3907.25
# Row 32
$ws.Cells.Item(32, 8).Value = 3577014.8
$ws.Cells.Item(32, 9).Value = 3707830
$ws.Cells.Item(32, 11).Value = 3707830
$ws.Cells.Item(32, 13).Value = -3707543
# Row 61
$ws.Cells.Item(61, 8).Value = 7804.65
$ws.Cells.Item(61, 9).Value = 2340.3333
$ws.Cells.Item(61, 10).Value = 16001.125
$ws.Cells.Item(61, 11).Value = 2340.3333
$ws.Cells.Item(61, 12).Value = 16001.125
$ws.Cells.Item(61, 13).Value = -2128.3333
$ws.Cells.Item(61, 14).Value = -16425.125
# Row 74
$ws.Cells.Item(74, 8).Value = 35058.195
$ws.Cells.Item(74, 9).Value = 54556.895
$ws.Cells.Item(74, 10).Value = 4185.25
$ws.Cells.Item(74, 11).Value = 54556.895
$ws.Cells.Item(74, 12).Value = 4185.25
$ws.Cells.Item(74, 13).Value = -53682.895
$ws.Cells.Item(74, 14).Value = -5933.25
# Row 77
$ws.Cells.Item(77, 8).Value = 35058.195
$ws.Cells.Item(77, 9).Value = 54556.895
$ws.Cells.Item(77, 10).Value = 4185.25
$ws.Cells.Item(77, 11).Value = 272784.475
$ws.Cells.Item(77, 12).Value = 20926.25
$ws.Cells.Item(77, 13).Value = -268416.475
$ws.Cells.Item(77, 14).Value = -29662.25
# Row 102
$ws.Cells.Item(102, 8).Value = 22731410
$ws.Cells.Item(102, 9).Value = 41670216
$ws.Cells.Item(102, 10).Value = 4839.9
$ws.Cells.Item(102, 11).Value = 41670216
$ws.Cells.Item(102, 12).Value = 4839.9
$ws.Cells.Item(102, 13).Value = -41668594
$ws.Cells.Item(102, 14).Value = -8083.9
# Row 116
$ws.Cells.Item(116, 8).Value = 5104.857
$ws.Cells.Item(116, 9).Value = 7003
$ws.Cells.Item(116, 10).Value = 3681.25
$ws.Cells.Item(116, 11).Value = 7003
$ws.Cells.Item(116, 12).Value = 3681.25
$ws.Cells.Item(116, 13).Value = -4709
$ws.Cells.Item(116, 14).Value = -8269.25
# Row 132
$ws.Cells.Item(132, 8).Value = 8052.205
$ws.Cells.Item(132, 9).Value = 6586.7407
$ws.Cells.Item(132, 10).Value = 11349.5
$ws.Cells.Item(132, 11).Value = 19760.2221
$ws.Cells.Item(132, 12).Value = 34048.5
$ws.Cells.Item(132, 13).Value = -17230.2221
$ws.Cells.Item(132, 14).Value = -39108.5
# Row 136
$ws.Cells.Item(136, 8).Value = 7804.65
$ws.Cells.Item(136, 9).Value = 2340.3333
$ws.Cells.Item(136, 10).Value = 16001.125
$ws.Cells.Item(136, 11).Value = 7020.999899999999
$ws.Cells.Item(136, 12).Value = 48003.375
$ws.Cells.Item(136, 13).Value = -4470.999899999999
$ws.Cells.Item(136, 14).Value = -53103.375

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 5104.857
$ws.Cells.Item(3, 9).Value = 7003
$ws.Cells.Item(3, 10).Value = 3681.25
$ws.Cells.Item(3, 11).Value = 7003
$ws.Cells.Item(3, 12).Value = 3681.25
$ws.Cells.Item(3, 13).Value = -6889
$ws.Cells.Item(3, 14).Value = -3909.25
# Row 36
$ws.Cells.Item(36, 8).Value = 466.33334
$ws.Cells.Item(36, 9).Value = 466.33334
$ws.Cells.Item(36, 11).Value = 466.33334
$ws.Cells.Item(36, 13).Value = 67.66665999999998
# Row 86
$ws.Cells.Item(86, 8).Value = 32863.812
$ws.Cells.Item(86, 9).Value = 50832.4
$ws.Cells.Item(86, 10).Value = 2916.1667
$ws.Cells.Item(86, 11).Value = 50832.4
$ws.Cells.Item(86, 12).Value = 2916.1667
$ws.Cells.Item(86, 13).Value = -49709.4
$ws.Cells.Item(86, 14).Value = -5162.1667
# Row 89
$ws.Cells.Item(89, 8).Value = 32863.812
$ws.Cells.Item(89, 9).Value = 50832.4
$ws.Cells.Item(89, 10).Value = 2916.1667
$ws.Cells.Item(89, 11).Value = 254162
$ws.Cells.Item(89, 12).Value = 14580.8335
$ws.Cells.Item(89, 13).Value = -248546
$ws.Cells.Item(89, 14).Value = -25812.8335
# Row 94
$ws.Cells.Item(94, 8).Value = 4848.3335
$ws.Cells.Item(94, 9).Value = 3375.7144
$ws.Cells.Item(94, 11).Value = 3375.7144
$ws.Cells.Item(94, 13).Value = -2924.7144
# Row 107
$ws.Cells.Item(107, 8).Value = 125012550
$ws.Cells.Item(107, 9).Value = 140638500
$ws.Cells.Item(107, 11).Value = 140638500
$ws.Cells.Item(107, 13).Value = -140636580
# Row 134
$ws.Cells.Item(134, 8).Value = 6068.108
$ws.Cells.Item(134, 9).Value = 2687.8262
$ws.Cells.Item(134, 11).Value = 8063.4786
$ws.Cells.Item(134, 13).Value = -5528.4786

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 6295.7144
$ws.Cells.Item(31, 9).Value = 1707.76
$ws.Cells.Item(31, 10).Value = 11074.833
$ws.Cells.Item(31, 11).Value = 1707.76
$ws.Cells.Item(31, 12).Value = 11074.833
$ws.Cells.Item(31, 13).Value = -1412.76
$ws.Cells.Item(31, 14).Value = -11664.833
# Row 34
$ws.Cells.Item(34, 8).Value = 6295.7144
$ws.Cells.Item(34, 9).Value = 1707.76
$ws.Cells.Item(34, 10).Value = 11074.833
$ws.Cells.Item(34, 11).Value = 1707.76
$ws.Cells.Item(34, 12).Value = 11074.833
$ws.Cells.Item(34, 13).Value = -1505.76
$ws.Cells.Item(34, 14).Value = -11478.833
# Row 122
$ws.Cells.Item(122, 8).Value = 2726.7273
$ws.Cells.Item(122, 10).Value = 3742
$ws.Cells.Item(122, 12).Value = 11226
$ws.Cells.Item(122, 14).Value = -16126
# Row 132
$ws.Cells.Item(132, 8).Value = 5991.0293
$ws.Cells.Item(132, 9).Value = 3133.4614
$ws.Cells.Item(132, 10).Value = 7760
$ws.Cells.Item(132, 11).Value = 9400.3842
$ws.Cells.Item(132, 12).Value = 23280
$ws.Cells.Item(132, 13).Value = -6870.3842
$ws.Cells.Item(132, 14).Value = -28340

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Cells.Item(122, 8).Value = 3144010.2
$ws.Cells.Item(122, 9).Value = 5658415.5
$ws.Cells.Item(122, 10).Value = 1003.75
$ws.Cells.Item(122, 11).Value = 50925739.5
$ws.Cells.Item(122, 12).Value = 9033.75
$ws.Cells.Item(122, 13).Value = -50923289.5
$ws.Cells.Item(122, 14).Value = -13933.75
# Row 131
$ws.Cells.Item(131, 8).Value = 3316.946
$ws.Cells.Item(131, 9).Value = 2250
$ws.Cells.Item(131, 10).Value = 3411.0881
$ws.Cells.Item(131, 11).Value = 6750
$ws.Cells.Item(131, 12).Value = 10233.2643
$ws.Cells.Item(131, 13).Value = -1710
$ws.Cells.Item(131, 14).Value = -20313.2643

$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Cells.Item(123, 8).Value = 29998.8
$ws.Cells.Item(123, 10).Value = 29998.8
$ws.Cells.Item(123, 12).Value = 29998.8
$ws.Cells.Item(123, 14).Value = -34898.8
# Row 132
$ws.Cells.Item(132, 8).Value = 5108.65
$ws.Cells.Item(132, 9).Value = 1828.9231
$ws.Cells.Item(132, 10).Value = 11199.571
$ws.Cells.Item(132, 11).Value = 5486.7693
$ws.Cells.Item(132, 12).Value = 33598.713
$ws.Cells.Item(132, 13).Value = -2956.7693
$ws.Cells.Item(132, 14).Value = -38658.713

$ws = $wb.Worksheets.Item("LTW")
# Row 107
$ws.Cells.Item(107, 8).Value = 3954
$ws.Cells.Item(107, 9).Value = 3954
$ws.Cells.Item(107, 11).Value = 3954
$ws.Cells.Item(107, 13).Value = -2034

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Cells.Item(113, 8).Value = 749.8929000000001
$ws.Cells.Item(113, 9).Value = 677.9737
$ws.Cells.Item(113, 11).Value = 2033.9211
$ws.Cells.Item(113, 13).Value = 136.0789
# Row 132
$ws.Cells.Item(132, 8).Value = 6446.12
$ws.Cells.Item(132, 9).Value = 10726.7
$ws.Cells.Item(132, 10).Value = 3592.4
$ws.Cells.Item(132, 11).Value = 32180.1
$ws.Cells.Item(132, 12).Value = 10777.2
$ws.Cells.Item(132, 13).Value = -29650.1
$ws.Cells.Item(132, 14).Value = -15837.2
